$d = $word.ActiveDocument

# --- Hunk 1: title paragraph " Api Rest" -> " Api " + "Rest" (split run) ---
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Api Rest")
if ($found) {
    $start = $rng.Start
    $end = $rng.End
    # "Api Rest" -> split after "Api " (4 chars)
    $r2 = $d.Range($start + 4, $end)
    $r2.Font.Bold = 0
    $r2.Font.Bold = 1
}

Write-Output "done hunk1"
